$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 2-13 with new TPM-derived values (columns E-T) ---
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 9.080941333333334
$ws.Cells.Item(2, 8).Value = 27.242824
$ws.Cells.Item(2, 9).Value = 0.3647727776818407
$ws.Cells.Item(2, 10).Value = 0.3647727776818407
$ws.Cells.Item(2, 11).Value = 2
$ws.Cells.Item(2, 12).Value = 0.6666666666666666
$ws.Cells.Item(2, 13).Value = 0.2262196666666666
$ws.Cells.Item(2, 14).Value = 0.6786589999999999
$ws.Cells.Item(2, 15).Value = 0.03145179203784564
$ws.Cells.Item(2, 16).Value = 0.03145179203784564
$ws.Cells.Item(2, 17).Value = 2.054287521446222
$ws.Cells.Item(2, 18).Value = 18.488587693016
$ws.Cells.Item(2, 19).Value = 0.01147275754471655
$ws.Cells.Item(2, 20).Value = 0.01147275754471656

$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 9.080941333333334
$ws.Cells.Item(3, 8).Value = 27.242824
$ws.Cells.Item(3, 9).Value = 0.3647727776818407
$ws.Cells.Item(3, 10).Value = 0.3647727776818407
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 6.93109
$ws.Cells.Item(3, 14).Value = 20.79327
$ws.Cells.Item(3, 15).Value = 0.9636438974901603
$ws.Cells.Item(3, 16).Value = 0.9636438974901604
$ws.Cells.Item(3, 17).Value = 62.94082166605334
$ws.Cells.Item(3, 18).Value = 566.4673949944801
$ws.Cells.Item(3, 19).Value = 0.3515110611836407
$ws.Cells.Item(3, 20).Value = 0.3515110611836408

$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 9.080941333333334
$ws.Cells.Item(4, 8).Value = 27.242824
$ws.Cells.Item(4, 9).Value = 0.3647727776818407
$ws.Cells.Item(4, 10).Value = 0.3647727776818407
$ws.Cells.Item(4, 11).Value = 1
$ws.Cells.Item(4, 12).Value = 0.3333333333333333
$ws.Cells.Item(4, 13).Value = 0.03527466666666667
$ws.Cells.Item(4, 14).Value = 0.105824
$ws.Cells.Item(4, 15).Value = 0.004904310471994002
$ws.Cells.Item(4, 16).Value = 0.004904310471994003
$ws.Cells.Item(4, 17).Value = 0.3203271785528889
$ws.Cells.Item(4, 18).Value = 2.882944606976
$ws.Cells.Item(4, 19).Value = 0.001788958953483391
$ws.Cells.Item(4, 20).Value = 0.001788958953483392

$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 2.958532
$ws.Cells.Item(5, 8).Value = 8.875596
$ws.Cells.Item(5, 9).Value = 0.1188414169728452
$ws.Cells.Item(5, 10).Value = 0.1188414169728452
$ws.Cells.Item(5, 11).Value = 2
$ws.Cells.Item(5, 12).Value = 0.6666666666666666
$ws.Cells.Item(5, 13).Value = 0.2262196666666666
$ws.Cells.Item(5, 14).Value = 0.6786589999999999
$ws.Cells.Item(5, 15).Value = 0.03145179203784564
$ws.Cells.Item(5, 16).Value = 0.03145179203784564
$ws.Cells.Item(5, 17).Value = 0.6692781228626665
$ws.Cells.Item(5, 18).Value = 6.023503105763999
$ws.Cells.Item(5, 19).Value = 0.003737775532112826
$ws.Cells.Item(5, 20).Value = 0.003737775532112827

$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 2.958532
$ws.Cells.Item(6, 8).Value = 8.875596
$ws.Cells.Item(6, 9).Value = 0.1188414169728452
$ws.Cells.Item(6, 10).Value = 0.1188414169728452
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 6.93109
$ws.Cells.Item(6, 14).Value = 20.79327
$ws.Cells.Item(6, 15).Value = 0.9636438974901603
$ws.Cells.Item(6, 16).Value = 0.9636438974901604
$ws.Cells.Item(6, 17).Value = 20.50585155988
$ws.Cells.Item(6, 18).Value = 184.55266403892
$ws.Cells.Item(6, 19).Value = 0.1145208062349658
$ws.Cells.Item(6, 20).Value = 0.1145208062349659

$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 2.958532
$ws.Cells.Item(7, 8).Value = 8.875596
$ws.Cells.Item(7, 9).Value = 0.1188414169728452
$ws.Cells.Item(7, 10).Value = 0.1188414169728452
$ws.Cells.Item(7, 11).Value = 1
$ws.Cells.Item(7, 12).Value = 0.3333333333333333
$ws.Cells.Item(7, 13).Value = 0.03527466666666667
$ws.Cells.Item(7, 14).Value = 0.105824
$ws.Cells.Item(7, 15).Value = 0.004904310471994002
$ws.Cells.Item(7, 16).Value = 0.004904310471994003
$ws.Cells.Item(7, 17).Value = 0.1043612301226667
$ws.Cells.Item(7, 18).Value = 0.939251071104
$ws.Cells.Item(7, 19).Value = 0.0005828352057665304
$ws.Cells.Item(7, 20).Value = 0.0005828352057665306

$ws.Cells.Item(8, 5).Value = 1
$ws.Cells.Item(8, 6).Value = 0.3333333333333333
$ws.Cells.Item(8, 7).Value = 0.1011933333333333
$ws.Cells.Item(8, 8).Value = 0.30358
$ws.Cells.Item(8, 9).Value = 0.004064839968450158
$ws.Cells.Item(8, 10).Value = 0.004064839968450158
$ws.Cells.Item(8, 11).Value = 2
$ws.Cells.Item(8, 12).Value = 0.6666666666666666
$ws.Cells.Item(8, 13).Value = 0.2262196666666666
$ws.Cells.Item(8, 14).Value = 0.6786589999999999
$ws.Cells.Item(8, 15).Value = 0.03145179203784564
$ws.Cells.Item(8, 16).Value = 0.03145179203784564
$ws.Cells.Item(8, 17).Value = 0.02289192213555555
$ws.Cells.Item(8, 18).Value = 0.20602729922
$ws.Cells.Item(8, 19).Value = 0.0001278465013548174
$ws.Cells.Item(8, 20).Value = 0.0001278465013548174

$ws.Cells.Item(9, 5).Value = 1
$ws.Cells.Item(9, 6).Value = 0.3333333333333333
$ws.Cells.Item(9, 7).Value = 0.1011933333333333
$ws.Cells.Item(9, 8).Value = 0.30358
$ws.Cells.Item(9, 9).Value = 0.004064839968450158
$ws.Cells.Item(9, 10).Value = 0.004064839968450158
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 6.93109
$ws.Cells.Item(9, 14).Value = 20.79327
$ws.Cells.Item(9, 15).Value = 0.9636438974901603
$ws.Cells.Item(9, 16).Value = 0.9636438974901604
$ws.Cells.Item(9, 17).Value = 0.7013801007333335
$ws.Cells.Item(9, 18).Value = 6.3124209066
$ws.Cells.Item(9, 19).Value = 0.00391705822987109
$ws.Cells.Item(9, 20).Value = 0.003917058229871091

$ws.Cells.Item(10, 5).Value = 1
$ws.Cells.Item(10, 6).Value = 0.3333333333333333
$ws.Cells.Item(10, 7).Value = 0.1011933333333333
$ws.Cells.Item(10, 8).Value = 0.30358
$ws.Cells.Item(10, 9).Value = 0.004064839968450158
$ws.Cells.Item(10, 10).Value = 0.004064839968450158
$ws.Cells.Item(10, 11).Value = 1
$ws.Cells.Item(10, 12).Value = 0.3333333333333333
$ws.Cells.Item(10, 13).Value = 0.03527466666666667
$ws.Cells.Item(10, 14).Value = 0.105824
$ws.Cells.Item(10, 15).Value = 0.004904310471994002
$ws.Cells.Item(10, 16).Value = 0.004904310471994003
$ws.Cells.Item(10, 17).Value = 0.003569561102222223
$ws.Cells.Item(10, 18).Value = 0.03212604992
$ws.Cells.Item(10, 19).Value = 0.00001993523722424988
$ws.Cells.Item(10, 20).Value = 0.00001993523722424989

$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 12.559037
$ws.Cells.Item(11, 8).Value = 37.677111
$ws.Cells.Item(11, 9).Value = 0.5044845730566345
$ws.Cells.Item(11, 10).Value = 0.5044845730566345
$ws.Cells.Item(11, 11).Value = 2
$ws.Cells.Item(11, 12).Value = 0.6666666666666666
$ws.Cells.Item(11, 13).Value = 0.2262196666666666
$ws.Cells.Item(11, 14).Value = 0.6786589999999999
$ws.Cells.Item(11, 15).Value = 0.03145179203784564
$ws.Cells.Item(11, 16).Value = 0.03145179203784564
$ws.Cells.Item(11, 17).Value = 2.841101163794333
$ws.Cells.Item(11, 18).Value = 25.569910474149
$ws.Cells.Item(11, 19).Value = 0.01586694387807861
$ws.Cells.Item(11, 20).Value = 0.01586694387807861

$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 12.559037
$ws.Cells.Item(12, 8).Value = 37.677111
$ws.Cells.Item(12, 9).Value = 0.5044845730566345
$ws.Cells.Item(12, 10).Value = 0.5044845730566345
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 6.93109
$ws.Cells.Item(12, 14).Value = 20.79327
$ws.Cells.Item(12, 15).Value = 0.9636438974901603
$ws.Cells.Item(12, 16).Value = 0.9636438974901604
$ws.Cells.Item(12, 17).Value = 87.04781576033001
$ws.Cells.Item(12, 18).Value = 783.4303418429701
$ws.Cells.Item(12, 19).Value = 0.4861434802039548
$ws.Cells.Item(12, 20).Value = 0.4861434802039549

$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 12.559037
$ws.Cells.Item(13, 8).Value = 37.677111
$ws.Cells.Item(13, 9).Value = 0.5044845730566345
$ws.Cells.Item(13, 10).Value = 0.5044845730566345
$ws.Cells.Item(13, 11).Value = 1
$ws.Cells.Item(13, 12).Value = 0.3333333333333333
$ws.Cells.Item(13, 13).Value = 0.03527466666666667
$ws.Cells.Item(13, 14).Value = 0.105824
$ws.Cells.Item(13, 15).Value = 0.004904310471994002
$ws.Cells.Item(13, 16).Value = 0.004904310471994003
$ws.Cells.Item(13, 17).Value = 0.4430158438293335
$ws.Cells.Item(13, 18).Value = 3.987142594464
$ws.Cells.Item(13, 19).Value = 0.002474148974601076
$ws.Cells.Item(13, 20).Value = 0.002474148974601076

# --- Add new rows 14-16 for Neutrophils sending cluster ---
$ws.Cells.Item(14, 1).Value = "Neutrophils"
$ws.Cells.Item(14, 2).Value = "Vtn"
$ws.Cells.Item(14, 3).Value = "Tnfrsf11b"
$ws.Cells.Item(14, 4).Value = "ECs"
$ws.Cells.Item(14, 5).Value = 1
$ws.Cells.Item(14, 6).Value = 0.3333333333333333
$ws.Cells.Item(14, 7).Value = 0.1950853333333333
$ws.Cells.Item(14, 8).Value = 0.585256
$ws.Cells.Item(14, 9).Value = 0.007836392320229479
$ws.Cells.Item(14, 10).Value = 0.007836392320229481
$ws.Cells.Item(14, 11).Value = 2
$ws.Cells.Item(14, 12).Value = 0.6666666666666666
$ws.Cells.Item(14, 13).Value = 0.2262196666666666
$ws.Cells.Item(14, 14).Value = 0.6786589999999999
$ws.Cells.Item(14, 15).Value = 0.03145179203784564
$ws.Cells.Item(14, 16).Value = 0.03145179203784564
$ws.Cells.Item(14, 17).Value = 0.04413213907822221
$ws.Cells.Item(14, 18).Value = 0.3971892517039999
$ws.Cells.Item(14, 19).Value = 0.0002464685815828283
$ws.Cells.Item(14, 20).Value = 0.0002464685815828283

$ws.Cells.Item(15, 1).Value = "Neutrophils"
$ws.Cells.Item(15, 2).Value = "Vtn"
$ws.Cells.Item(15, 3).Value = "Tnfrsf11b"
$ws.Cells.Item(15, 4).Value = "FAPs"
$ws.Cells.Item(15, 5).Value = 1
$ws.Cells.Item(15, 6).Value = 0.3333333333333333
$ws.Cells.Item(15, 7).Value = 0.1950853333333333
$ws.Cells.Item(15, 8).Value = 0.585256
$ws.Cells.Item(15, 9).Value = 0.007836392320229479
$ws.Cells.Item(15, 10).Value = 0.007836392320229481
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 12).Value = 1
$ws.Cells.Item(15, 13).Value = 6.93109
$ws.Cells.Item(15, 14).Value = 20.79327
$ws.Cells.Item(15, 15).Value = 0.9636438974901603
$ws.Cells.Item(15, 16).Value = 0.9636438974901604
$ws.Cells.Item(15, 17).Value = 1.352154003013333
$ws.Cells.Item(15, 18).Value = 12.16938602712
$ws.Cells.Item(15, 19).Value = 0.007551491637727896
$ws.Cells.Item(15, 20).Value = 0.007551491637727899

$ws.Cells.Item(16, 1).Value = "Neutrophils"
$ws.Cells.Item(16, 2).Value = "Vtn"
$ws.Cells.Item(16, 3).Value = "Tnfrsf11b"
$ws.Cells.Item(16, 4).Value = "MuSCs"
$ws.Cells.Item(16, 5).Value = 1
$ws.Cells.Item(16, 6).Value = 0.3333333333333333
$ws.Cells.Item(16, 7).Value = 0.1950853333333333
$ws.Cells.Item(16, 8).Value = 0.585256
$ws.Cells.Item(16, 9).Value = 0.007836392320229479
$ws.Cells.Item(16, 10).Value = 0.007836392320229481
$ws.Cells.Item(16, 11).Value = 1
$ws.Cells.Item(16, 12).Value = 0.3333333333333333
$ws.Cells.Item(16, 13).Value = 0.03527466666666667
$ws.Cells.Item(16, 14).Value = 0.105824
$ws.Cells.Item(16, 15).Value = 0.004904310471994002
$ws.Cells.Item(16, 16).Value = 0.004904310471994003
$ws.Cells.Item(16, 17).Value = 0.006881570104888889
$ws.Cells.Item(16, 18).Value = 0.061934130944
$ws.Cells.Item(16, 19).Value = 0.00003843210091875481
$ws.Cells.Item(16, 20).Value = 0.00003843210091875483

Write-Output "edit complete"